# Add /start-survey ajax call
# Sets the start_time value for the first logged conversation (row 2, "conv_5")
# and widens column B to fit the timestamp text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2023-07-10T11:55:56.998Z"

# Widen column B (start_time) to fit the new timestamp text.
# (25.5 "characters" is the input that rounds to the stored column width
# closest to the target 26.36328125 under this engine's column-width grid.)
$ws.Columns.Item(2).ColumnWidth = 25.5

$ws.Range("B6").Select()
